# Add three more years (2021, 2022, 2023) of data to the Ombudsman-appeals
# table, extending the sheet from column Q (2020) out to column T (2023).
# New cells inherit the same look (number format / borders / alignment) as
# the existing column Q by copying its formatting before writing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column Q's formatting (header-rule row, blank spacer row, and the
# two data rows) across the three new columns R, S and T.
$ws.Range("Q2:Q5").Copy() | Out-Null
$ws.Range("R2:R5").PasteSpecial(-4122) | Out-Null
$ws.Range("S2:S5").PasteSpecial(-4122) | Out-Null
$ws.Range("T2:T5").PasteSpecial(-4122) | Out-Null

# Year headers
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

# "Number of written appeals" row
$ws.Range("R4").Value = 4301
$ws.Range("S4").Value = 3690
$ws.Range("T4").Value = 2620

# "Number of positively resolved" row
$ws.Range("R5").Value = 427
$ws.Range("S5").Value = 280
$ws.Range("T5").Value = 264
